# Add a 4pt ("before") space-before setting to the "Body Text" paragraph
# style (w:spacing w:before="80", 80 twentieths-of-a-point == 4pt).
#
# Commit message:
#   Add space-before to text-box generated from floatfig
#   Also define 4pt as space-before to "Body Text" style.
#
# This document has no floating-figure textbox content, so the only
# observable, scriptable change here is the "Body Text" style update.

$d = $word.ActiveDocument

$bodyText = $d.Styles("Body Text")
$bodyText.ParagraphFormat.SpaceBefore = 4
